# The sheet is a weekly price log for "Poroto verde" at Macroferia Regional
# de Talca. This commit adds one new weekly record. It is inserted as a new
# row 69 (pushing the existing rows 69-134 down to 70-135), which is why the
# diff shows every row from 69 onward "shifting" to the values that used to
# belong to the row above it, and a brand new row 135 appears holding what
# used to be the last row (134)'s data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 69; existing rows 69:134 shift down to 70:135,
# inheriting formatting (incl. the date style on column D) from the row above.
$ws.Rows("69").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A69").Value = 5
$ws.Range("B69").Value = "Macroferia Regional de Talca"
$ws.Range("C69").Value = "Maule"
$ws.Range("D69").Value = 44566
$ws.Range("E69").Value = 7
$ws.Range("F69").Value = 100112031
$ws.Range("G69").Value = "Poroto verde"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 200
$ws.Range("K69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("M69").Value = 30000
$ws.Range("N69").Value = "`$/saco 25 kilos"
$ws.Range("O69").Value = "Región del Maule"
$ws.Range("P69").Value = 1200
$ws.Range("Q69").Value = 25
$ws.Range("R69").Value = "Hortaliza"
